$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential disclosure date from 2021-05-20 to 2021-05-21
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-21 for illustrative purposes only and are subject to change."

# Update the Weight / Percent Change figures
$ws.Range("D2").Value = 0.8479035817348498
$ws.Range("E2").Value = 0.001376376376376554

$ws.Range("D3").Value = 0.1520964182651502
$ws.Range("E3").Value = -0.01153273809523803

$ws.Range("E4").Value = -0.0005870536977309326

# Restore sheet protection (unprotected above to allow the edits)
$ws.Protect()
